# Refresh market-price-derived profit columns (H-N) across all job sheets.
# Source values come from a scheduled price-lookup run; no formulas are
# involved anywhere in this workbook, so every cell below is a literal write.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 70: Consecrating Congregation / Holy Water
$ws.Range("H70").Value = 3202.2222
$ws.Range("J70").Value = 2579.2856
$ws.Range("L70").Value = 7737.8568
$ws.Range("N70").Value = -8277.856800000001
# Row 73: Curbing the Contagion (L) / Holy Water
$ws.Range("H73").Value = 3202.2222
$ws.Range("J73").Value = 2579.2856
$ws.Range("L73").Value = 7737.8568
$ws.Range("N73").Value = -9609.856800000001
# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 4179.4614
$ws.Range("I88").Value = 5499.8
$ws.Range("J88").Value = 3354.25
$ws.Range("K88").Value = 5499.8
$ws.Range("L88").Value = 3354.25
$ws.Range("M88").Value = -5093.8
$ws.Range("N88").Value = -4166.25
# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 4179.4614
$ws.Range("I91").Value = 5499.8
$ws.Range("J91").Value = 3354.25
$ws.Range("K91").Value = 5499.8
$ws.Range("L91").Value = 3354.25
$ws.Range("M91").Value = -4095.8
$ws.Range("N91").Value = -6162.25
# Row 92: Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 1959.1875
$ws.Range("I92").Value = 2054.3
$ws.Range("K92").Value = 2054.3
$ws.Range("M92").Value = -806.3000000000002
# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 4569.9
$ws.Range("J112").Value = 4855.4443
$ws.Range("L112").Value = 14566.3329
$ws.Range("N112").Value = -16782.3329
# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 10416.429
$ws.Range("J113").Value = 5459.875
$ws.Range("L113").Value = 5459.875
$ws.Range("N113").Value = -11967.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 3840.85
$ws.Range("J2").Value = 2776.625
$ws.Range("L2").Value = 2776.625
$ws.Range("N2").Value = -3002.625
# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 3840.85
$ws.Range("J116").Value = 2776.625
$ws.Range("L116").Value = 2776.625
$ws.Range("N116").Value = -7364.625
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 3718620.8
$ws.Range("I132").Value = 4453.7617
$ws.Range("K132").Value = 13361.2851
$ws.Range("M132").Value = -10831.2851

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 3840.85
$ws.Range("J3").Value = 2776.625
$ws.Range("L3").Value = 2776.625
$ws.Range("N3").Value = -3004.625
# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 22960.074
$ws.Range("I99").Value = 20196.88
$ws.Range("K99").Value = 20196.88
$ws.Range("M99").Value = -18698.88

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 20: Re-crating the Scene / Iron Spear
$ws.Range("H20").Value = 101000
$ws.Range("J20").Value = 101000
$ws.Range("L20").Value = 101000
$ws.Range("N20").Value = -101472
# Row 30: Polearms Aplenty / Iron Spear
$ws.Range("H30").Value = 101000
$ws.Range("J30").Value = 101000
$ws.Range("L30").Value = 101000
$ws.Range("N30").Value = -101182
# Row 48: The Cold, Cold Ground / Oak Pattens
$ws.Range("H48").Value = 46842
$ws.Range("J48").Value = 46842
$ws.Range("L48").Value = 46842
$ws.Range("N48").Value = -47794
# Row 76: Walking on Pins and Needles / Dark Chestnut Lumber
$ws.Range("H76").Value = 8187.5
$ws.Range("I76").Value = 8187.5
$ws.Range("K76").Value = 8187.5
$ws.Range("M76").Value = -7872.5
# Row 79: Like Lemon on a Lumbercut (L) / Dark Chestnut Lumber
$ws.Range("H79").Value = 8187.5
$ws.Range("I79").Value = 8187.5
$ws.Range("K79").Value = 8187.5
$ws.Range("M79").Value = -7095.5
# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 3453566.5
$ws.Range("I99").Value = 2315408
$ws.Range("J99").Value = 6678348.5
$ws.Range("K99").Value = 2315408
$ws.Range("L99").Value = 6678348.5
$ws.Range("M99").Value = -2313910
$ws.Range("N99").Value = -6681344.5
# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 3453566.5
$ws.Range("I126").Value = 2315408
$ws.Range("J126").Value = 6678348.5
$ws.Range("K126").Value = 6946224
$ws.Range("L126").Value = 20035045.5
$ws.Range("M126").Value = -6943754
$ws.Range("N126").Value = -20039985.5
# Row 128: An A-prop-riate Request / Ironwood Spear
$ws.Range("H128").Value = 101000
$ws.Range("J128").Value = 101000
$ws.Range("L128").Value = 101000
$ws.Range("N128").Value = -110960
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2661
$ws.Range("I132").Value = 2661
$ws.Range("K132").Value = 7983
$ws.Range("M132").Value = -5453

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 104: Fits to a Tea / Doman Tea
$ws.Range("H104").Value = 6174339.5
$ws.Range("I104").Value = 3000
$ws.Range("J104").Value = 12345679
$ws.Range("K104").Value = 9000
$ws.Range("L104").Value = 37037037
$ws.Range("M104").Value = -6379
$ws.Range("N104").Value = -37042279
# Row 117: A Good Omen / Peppered Popotoes
$ws.Range("H117").Value = 514.6667
$ws.Range("I117").Value = 514.6667
$ws.Range("K117").Value = 1544.0001
$ws.Range("M117").Value = 1897.9999
# Row 129: Comfort Food / Yakow Moussaka
$ws.Range("H129").Value = 12988132
$ws.Range("I129").Value = 1366.8
$ws.Range("J129").Value = 45455044
$ws.Range("K129").Value = 4100.4
$ws.Range("L129").Value = 136365132
$ws.Range("M129").Value = 899.6000000000004
$ws.Range("N129").Value = -136375132
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1464.06
$ws.Range("I131").Value = 1076.6666
$ws.Range("J131").Value = 1476.0413
$ws.Range("K131").Value = 3229.9998
$ws.Range("L131").Value = 4428.123900000001
$ws.Range("M131").Value = 1810.0002
$ws.Range("N131").Value = -14508.1239

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 21220.9
$ws.Range("J70").Value = 21003
$ws.Range("L70").Value = 21003
$ws.Range("N70").Value = -21543
# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 21220.9
$ws.Range("J73").Value = 21003
$ws.Range("L73").Value = 21003
$ws.Range("N73").Value = -22875
# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 1056.2727
$ws.Range("I97").Value = 998.9259
$ws.Range("K97").Value = 998.9259
$ws.Range("M97").Value = -502.9259
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 4035574.8
$ws.Range("I102").Value = 7817246
$ws.Range("J102").Value = 1791.8667
$ws.Range("K102").Value = 7817246
$ws.Range("L102").Value = 1791.8667
$ws.Range("M102").Value = -7815624
$ws.Range("N102").Value = -5035.8667
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 4674597.5
$ws.Range("I126").Value = 3601257.2
$ws.Range("K126").Value = 10803771.6
$ws.Range("M126").Value = -10801301.6

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 23811462
$ws.Range("I22").Value = 1472.591
$ws.Range("J22").Value = 50002450
$ws.Range("K22").Value = 1472.591
$ws.Range("L22").Value = 50002450
$ws.Range("M22").Value = -1177.591
$ws.Range("N22").Value = -50003040
# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 23811462
$ws.Range("I27").Value = 1472.591
$ws.Range("J27").Value = 50002450
$ws.Range("K27").Value = 1472.591
$ws.Range("L27").Value = 50002450
$ws.Range("M27").Value = -1365.591
$ws.Range("N27").Value = -50002664
# Row 50: The Birdmen of Ishgard / Boarskin Culottes
$ws.Range("H50").Value = 60000
$ws.Range("J50").Value = 60000
$ws.Range("L50").Value = 60000
$ws.Range("N50").Value = -61274
# Row 54: Not So Alike in Dignity / Boarskin Jerkin
$ws.Range("H54").Value = 20420
$ws.Range("J54").Value = 20420
$ws.Range("L54").Value = 20420
$ws.Range("N54").Value = -21708
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 23228270
$ws.Range("I122").Value = 43399764
$ws.Range("K122").Value = 130199292
$ws.Range("M122").Value = -130196842
# Row 138: Freezing Toes / Gomphotherium Boots of Striking
$ws.Range("H138").Value = 49999
$ws.Range("J138").Value = 49999
$ws.Range("L138").Value = 49999
$ws.Range("N138").Value = -60279

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 133: Begin with the Basics / Snow Cotton Jacket
$ws.Range("H133").Value = 54783.832
$ws.Range("J133").Value = 54783.832
$ws.Range("L133").Value = 54783.832
$ws.Range("N133").Value = -64903.832
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 16508.062
$ws.Range("I136").Value = 4493.2915
$ws.Range("J136").Value = 52552.375
$ws.Range("K136").Value = 13479.8745
$ws.Range("L136").Value = 157657.125
$ws.Range("M136").Value = -10929.8745
$ws.Range("N136").Value = -162757.125
